$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "07.11.2024  (*07:00*)"
$ws.Range("F2").Value = "03.10 *08:15 - 10:00* (**LAB2**)"
$ws.Range("F3").Value = "10.10 *08:15 - 10:00* (**LAB2**)"
$ws.Range("F4").Value = "17.10 *08:15 - 10:00* (**AUD G**)"
$ws.Range("F5").Value = "24.10 *08:15 - 10:00* (**LAB2**)"
$ws.Range("F6").Value = "31.10 *08:15 - 10:00* (**LAB2**)"

$ws.Range("F7").Style = "Normal"
$ws.Range("F7").Value = "07:11 *08:15 - 10:00* (**LAB2**)"

$ws.Range("F8").Select()
